$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 542.2143
$ws.Range("I39").Value = 169
$ws.Range("J39").Value = 822.125
$ws.Range("K39").Value = 507
$ws.Range("L39").Value = 2466.375
$ws.Range("M39").Value = -211
$ws.Range("N39").Value = -3058.375
$ws.Range("H57").Value = 24339.75
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 24339.75
$ws.Range("K57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("M57").Value = 73019.25
$ws.Range("N57").Value = -74017.25
$ws.Range("H129").Value = 997.97015
$ws.Range("I129").Value = 1023.1111
$ws.Range("J129").Value = 994.069
$ws.Range("K129").Value = 3069.3333
$ws.Range("L129").Value = 2982.207
$ws.Range("M129").Value = 1930.6667
$ws.Range("N129").Value = -12982.207
$ws.Range("H131").Value = 5436.8184
$ws.Range("I131").Value = 947.5
$ws.Range("J131").Value = 6434.4443
$ws.Range("K131").Value = 2842.5
$ws.Range("L131").Value = 19303.3329
$ws.Range("M131").Value = 2197.5
$ws.Range("N131").Value = -29383.3329
$ws.Range("H138").Value = 2257.4038
$ws.Range("I138").Value = 1216.625
$ws.Range("K138").Value = 3649.875
$ws.Range("M138").Value = 1490.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 112873.11
$ws.Range("I45").Value = 168091.33
$ws.Range("J45").Value = 2436.6667
$ws.Range("K45").Value = 168091.33
$ws.Range("L45").Value = 2436.6667
$ws.Range("M45").Value = -167714.33
$ws.Range("N45").Value = -3190.6667
$ws.Range("H54").Value = 21033.334
$ws.Range("I54").Value = 28500
$ws.Range("J54").Value = 20500
$ws.Range("K54").Value = 28500
$ws.Range("L54").Value = 20500
$ws.Range("M54").Value = -27731
$ws.Range("N54").Value = -22038
$ws.Range("H74").Value = 1097.1526
$ws.Range("I74").Value = 1079.6586
$ws.Range("J74").Value = 1137
$ws.Range("K74").Value = 1079.6586
$ws.Range("L74").Value = 1137
$ws.Range("M74").Value = -205.6586
$ws.Range("N74").Value = -2885
$ws.Range("H77").Value = 1097.1526
$ws.Range("I77").Value = 1079.6586
$ws.Range("J77").Value = 1137
$ws.Range("K77").Value = 5398.293
$ws.Range("L77").Value = 5685
$ws.Range("M77").Value = -1030.293
$ws.Range("N77").Value = -14421
$ws.Range("H95").Value = 19169.334
$ws.Range("J95").Value = 19169.334
$ws.Range("L95").Value = 19169.334
$ws.Range("N95").Value = -24661.334
$ws.Range("H139").Value = 48357.5
$ws.Range("J139").Value = 48357.5
$ws.Range("L139").Value = 48357.5
$ws.Range("N139").Value = -58637.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 12682.857
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 12682.857
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = 12682.857
$ws.Range("N81").Value = -14804.857
$ws.Range("H84").Value = 12682.857
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 12682.857
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = 38048.571
$ws.Range("N84").Value = -48656.571
$ws.Range("H135").Value = 43840
$ws.Range("J135").Value = 43840
$ws.Range("L135").Value = 43840
$ws.Range("N135").Value = -53980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1463249.8
$ws.Range("I31").Value = 1852699.6
$ws.Range("K31").Value = 1852699.6
$ws.Range("M31").Value = -1852404.6
$ws.Range("H34").Value = 1463249.8
$ws.Range("I34").Value = 1852699.6
$ws.Range("K34").Value = 1852699.6
$ws.Range("M34").Value = -1852497.6
$ws.Range("H58").Value = 34483396
$ws.Range("I58").Value = 45455176
$ws.Range("K58").Value = 45455176
$ws.Range("M58").Value = -45454973
$ws.Range("H136").Value = 34483396
$ws.Range("I136").Value = 45455176
$ws.Range("K136").Value = 136365528
$ws.Range("M136").Value = -136362978

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 207.14285
$ws.Range("I17").Value = 207.14285
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 621.4285500000001
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -452.4285500000001
$ws.Range("H49").Value = 1822.2222
$ws.Range("J49").Value = 1822.2222
$ws.Range("L49").Value = 5466.6666
$ws.Range("N49").Value = -5778.6666
$ws.Range("H61").Value = 367.57144
$ws.Range("I61").Value = 44.5
$ws.Range("J61").Value = 798.3333
$ws.Range("K61").Value = 133.5
$ws.Range("L61").Value = 2394.9999
$ws.Range("M61").Value = 81.5
$ws.Range("N61").Value = -2824.9999
$ws.Range("H69").Value = 37500.332
$ws.Range("I69").Value = 390
$ws.Range("J69").Value = 56055.5
$ws.Range("K69").Value = 1170
$ws.Range("L69").Value = 168166.5
$ws.Range("M69").Value = -359
$ws.Range("N69").Value = -169788.5
$ws.Range("H72").Value = 37500.332
$ws.Range("I72").Value = 390
$ws.Range("J72").Value = 56055.5
$ws.Range("K72").Value = 3510
$ws.Range("L72").Value = 504499.5
$ws.Range("M72").Value = 546
$ws.Range("N72").Value = -512611.5
$ws.Range("H131").Value = 938.3099999999999
$ws.Range("I131").Value = 615
$ws.Range("J131").Value = 951.78125
$ws.Range("K131").Value = 1845
$ws.Range("L131").Value = 2855.34375
$ws.Range("M131").Value = 3195
$ws.Range("N131").Value = -12935.34375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 20093.143
$ws.Range("J134").Value = 20093.143
$ws.Range("L134").Value = 60279.429
$ws.Range("N134").Value = -65349.429
$ws.Range("H136").Value = 20839.715
$ws.Range("J136").Value = 20839.715
$ws.Range("L136").Value = 62519.145
$ws.Range("N136").Value = -67619.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3585.7144
$ws.Range("I22").Value = 3666.6667
$ws.Range("J22").Value = 3525
$ws.Range("K22").Value = 3666.6667
$ws.Range("L22").Value = 3525
$ws.Range("M22").Value = -3371.6667
$ws.Range("N22").Value = -4115
$ws.Range("H27").Value = 3585.7144
$ws.Range("I27").Value = 3666.6667
$ws.Range("J27").Value = 3525
$ws.Range("K27").Value = 3666.6667
$ws.Range("L27").Value = 3525
$ws.Range("M27").Value = -3559.6667
$ws.Range("N27").Value = -3739
$ws.Range("H132").Value = 45724756
$ws.Range("I132").Value = 103897260
$ws.Range("J132").Value = 17784.785
$ws.Range("K132").Value = 311691780
$ws.Range("L132").Value = 53354.355
$ws.Range("M132").Value = -311689250
$ws.Range("N132").Value = -58414.355
$ws.Range("H136").Value = 56393320
$ws.Range("I136").Value = 42329900
$ws.Range("J136").Value = 90912616
$ws.Range("K136").Value = 126989700
$ws.Range("L136").Value = 272737848
$ws.Range("M136").Value = -126987150
$ws.Range("N136").Value = -272742948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 47800.332
$ws.Range("J46").Value = 47800.332
$ws.Range("L46").Value = 47800.332
$ws.Range("N46").Value = -48262.332
$ws.Range("H107").Value = 1073.1666
$ws.Range("I107").Value = 756.4
$ws.Range("J107").Value = 1299.4286
$ws.Range("K107").Value = 2269.2
$ws.Range("L107").Value = 3898.2858
$ws.Range("M107").Value = -349.1999999999998
$ws.Range("N107").Value = -7738.2858
$ws.Range("H132").Value = 58485.9
$ws.Range("I132").Value = 80270.766
$ws.Range("J132").Value = 18028.285
$ws.Range("K132").Value = 240812.298
$ws.Range("L132").Value = 54084.855
$ws.Range("M132").Value = -238282.298
$ws.Range("N132").Value = -59144.855
$ws.Range("H134").Value = 47800.332
$ws.Range("J134").Value = 47800.332
$ws.Range("L134").Value = 143400.996
$ws.Range("N134").Value = -148470.996
